$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Detail Date" timestamp in column B for all data rows (2-38)
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 2).Value = "2023-06-27 22:30:24"
}

# Update iOS Rank values in column G (duplicates fix produced new rank numbers)
$ws.Cells.Item(5, 7).Value = 96
$ws.Cells.Item(7, 7).Value = 89
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(10, 7).Value = 80
$ws.Cells.Item(11, 7).Value = 66
$ws.Cells.Item(14, 7).Value = 104
$ws.Cells.Item(15, 7).Value = 150
$ws.Cells.Item(16, 7).Value = 130
$ws.Cells.Item(19, 7).Value = 12
$ws.Cells.Item(20, 7).Value = 141
$ws.Cells.Item(22, 7).Value = 57
$ws.Cells.Item(23, 7).Value = 184
